# unify the conception of DataNode, DataTable, Entity.
# Rename the sheet "Property1" -> "DataNode" to match the unified naming.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Leave the cursor where the author last left it when saving.
$ws.Range("D36").Select()
